# Updated symbol list on Sun Feb  5 21:35:52 UTC 2023 with GitHub Actions
# Refresh Price (D) and Volume(1h) (E) columns for the crypto rows that
# changed in this run. Values are plain-text cells (no numeric/percent
# typing), matching the source data feed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A cell whose style is the plain "no special formatting" default, used
# to strip the automatic text-number-format style Excel likes to stamp
# on cells once their content has to be forced to Text.
$plainStyleSource = $ws.Range("B2")

function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = $plainStyleSource.Style
}

# Map of row -> @{ D = newPrice; E = newVolume }
$updates = @{
    2  = @{ D = "326.93"; E = "-0.97%" }
    3  = @{ D = "43.69"; E = "4.89%" }
    4  = @{ E = "-3.56%" }
    5  = @{ D = "0.08066"; E = "-4.26%" }
    6  = @{ D = "8.633"; E = "-1.80%" }
    7  = @{ D = "4.280"; E = "-4.64%" }
    8  = @{ D = "1.876"; E = "-5.64%" }
    10 = @{ D = "0.9356"; E = "0.92%" }
    11 = @{ D = "0.1171"; E = "-8.24%" }
    12 = @{ D = "0.1895"; E = "-4.16%" }
    13 = @{ D = "0.09580"; E = "1.29%" }
    14 = @{ D = "0.04147"; E = "5.68%" }
    15 = @{ D = "0.1065"; E = "0.29%" }
    16 = @{ E = "-2.15%" }
    17 = @{ D = "0.005972"; E = "-2.34%" }
    18 = @{ D = "3.565"; E = "4.19%" }
    19 = @{ E = "-0.42%" }
    20 = @{ D = "8.551"; E = "-5.14%" }
    21 = @{ D = "0.1365"; E = "0.10%" }
    22 = @{ D = "0.2592"; E = "3.14%" }
    23 = @{ D = "0.04343"; E = "-1.63%" }
    24 = @{ D = "0.001233"; E = "-1.04%" }
    25 = @{ D = "0.004345"; E = "-0.74%" }
    26 = @{ D = "0.0001232"; E = "3.30%" }
    27 = @{ D = "0.0003996"; E = "0.00%" }
    39 = @{ D = "0.02659"; E = "-6.20%" }
    40 = @{ D = "0.05430"; E = "-1.67%" }
    41 = @{ D = "0.01144"; E = "27.33%" }
    42 = @{ D = "0.007677"; E = "-3.11%" }
    43 = @{ E = "-3.45%" }
    44 = @{ D = "0.002126"; E = "1.98%" }
    45 = @{ D = "0.009675"; E = "-12.05%" }
    46 = @{ D = "0.00006867"; E = "-5.63%" }
    47 = @{ D = "0.00000000751"; E = "0.01%" }
    48 = @{ D = "0.003567"; E = "9.72%" }
    49 = @{ D = "0.002274"; E = "-0.33%" }
    50 = @{ D = "0.00002104"; E = "0.01%" }
    51 = @{ D = "0.0002004"; E = "0.01%" }
}

foreach ($row in $updates.Keys) {
    $rowUpdates = $updates[$row]
    if ($rowUpdates.ContainsKey("D")) {
        Set-TextValue $ws.Range("D$row") $rowUpdates["D"]
    }
    if ($rowUpdates.ContainsKey("E")) {
        Set-TextValue $ws.Range("E$row") $rowUpdates["E"]
    }
}
